# Fix typo in the "800 - 900" bucket label (was missing a space before 900)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A7").Value = " 800 - 900"
